# "Results from R script" - appends the two newest OHLC rows (2024-07-01
# and 2024-07-02) produced by the data-refresh script to the ELSA.MI price
# history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# --- Row 193 -------------------------------------------------------------
$ws.Cells.Item(192, 1).Copy()
$ws.Cells.Item(193, 1).PasteSpecial(-4122)        # xlPasteFormats (date style)
$ws.Cells.Item(193, 1).Value = 45474.2916666667
$ws.Cells.Item(193, 2).Value = 0
$ws.Cells.Item(193, 3).Value = 3.74000000953674
$ws.Cells.Item(193, 4).Value = 3.74000000953674
$ws.Cells.Item(193, 5).Value = 3.74000000953674
$ws.Cells.Item(193, 6).Value = 3.74000000953674
$ws.Cells.Item(193, 8).Value = "ELSA.MI"

# --- Row 194 -------------------------------------------------------------
$ws.Cells.Item(192, 1).Copy()
$ws.Cells.Item(194, 1).PasteSpecial(-4122)        # xlPasteFormats (date style)
$ws.Cells.Item(194, 1).Value = 45475.6457060185
$ws.Cells.Item(194, 2).Value = 2500
$ws.Cells.Item(194, 3).Value = 3.96000003814697
$ws.Cells.Item(194, 4).Value = 3.77999997138977
$ws.Cells.Item(194, 5).Value = 3.77999997138977
$ws.Cells.Item(194, 6).Value = 3.90000009536743
$ws.Cells.Item(194, 8).Value = "ELSA.MI"

# --- adj_close (column G) values -----------------------------------------
# These columns hold numeric-looking values that are nevertheless stored as
# plain text/shared strings in the workbook. A scratch cell is used to coerce
# the value to text (via a Text number format) and only the resulting value
# is pasted onto the target cells, so the target cells keep their original
# (default/general) style - just like the existing rows above them.
$scratch = $ws.Cells.Item(500, 50)

$scratch.NumberFormat = "@"
$scratch.Value = "3.74000000953674"
$scratch.Copy()
$ws.Cells.Item(193, 7).PasteSpecial(-4163)        # xlPasteValues

$scratch.Value = "3.90000009536743"
$scratch.Copy()
$ws.Cells.Item(194, 7).PasteSpecial(-4163)        # xlPasteValues

$scratch.Clear()
